# Applies the "Updates Small Aperture Score & Fix a bug for gantry speed"
# commit to QA list.xlsx:
#   - adds a new "VMAT0010" column (C) value next to the existing
#     "VMAT0028" column (B) header row
#   - appends two new rows (49, 50) with "Max Leaf Speed" and
#     "MLC speed distribution" labels, styled like the other red labels
#   - widens the new column C
#   - updates the active selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header value in row 2, column C (VMAT0010), alongside existing
# VMAT0028 in column B.
$ws.Range("C2").Value = "VMAT0010"

# New column C width (matches the width Excel stores for the new column).
$ws.Range("C1:C50").ColumnWidth = 13.93

# Two new labeled rows appended after row 47 (row 48 intentionally left
# blank, matching the source workbook). Shared-string table order follows
# write order, so add "MLC speed distribution" (idx 48) before
# "Max Leaf Speed" (idx 49).
$ws.Range("A50").Value = "MLC speed distribution"
$ws.Range("A49").Value = "Max Leaf Speed"

# Copy the red "section label" style used by the other A-column labels
# (e.g. A47) onto the two new cells.
$ws.Range("A47").Copy()
$ws.Range("A49:A50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the view: active cell/selection then scroll position.
$ws.Range("B35").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
